$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")
$ws.Activate()

# Shorten the "Goto the Doctor..." text first (still at its original row 16).
$ws.Range("A16").Value2 = "Goto the Doctor and get checked up"

# --- Insert two new rows above row 16 ("Goto the Doctor...") ---
# These become rows 16 & 17 ("Cancel LogMeIn.com" and "Bank of America Maintenance Fees"),
# pushing "Goto the Doctor..." down to row 18.
$ws.Rows("16:17").Insert()

# Copy formatting from the (now shifted) "Goto the Doctor..." row (row 18) into the two new rows.
$ws.Range("A18:B18").Copy()
$ws.Range("A16:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the values for the two new rows.
$ws.Range("A17").Value2 = "Bank of America Maintenance Fees"
$ws.Range("B17").Value2 = "TODO"
$ws.Range("A16").Value2 = "Cancel LogMeIn.com"
$ws.Range("B16").Value2 = "TODO"

# --- Insert one new row after row 20 ("Honda Accord Jack...") for the MRI entry ---
$ws.Rows("21:21").Insert()
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A21").Value2 = "Get MRI Images from Kachar"
$ws.Range("B21").Value2 = "TODO"

# --- Column width updates ---
$ws.Columns("A").ColumnWidth = 61.6
$ws.Columns("B").ColumnWidth = 18.6

# --- Selection / view ---
$ws.Range("A21").Select()

Write-Output "done"
